$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.739.88"
$ws.Range("E2").Value = "  +9.19%  "
$ws.Range("D3").Value = "2.681.59"
$ws.Range("E3").Value = "  +10.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.78"
$ws.Range("E5").Value = "  +13.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "588.20"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +4.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.196"
$ws.Range("E9").Value = "  +16.08%  "
$ws.Range("D10").Value = "2.683.86"
$ws.Range("E10").Value = "  +10.66%  "
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  +7.41%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "75.476.82"
$ws.Range("E14").Value = "  +9.21%  "
$ws.Range("D15").Value = "3.170.01"
$ws.Range("E15").Value = "  +10.61%  "
$ws.Range("E16").Value = "  +5.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.51"
$ws.Range("E17").Value = "  +11.06%  "
$ws.Range("D18").Value = "2.690.75"
$ws.Range("E18").Value = "  +11.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.29"
$ws.Range("E19").Value = "  +31.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.97"
$ws.Range("E20").Value = "  +11.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.05"
$ws.Range("E21").Value = "  +9.55%  "
$ws.Range("E22").Value = "  +15.86%  "
$ws.Range("E23").Value = "  +5.45%  "
$ws.Range("E24").Value = "  +4.65%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.83"
$ws.Range("E26").Value = "  +6.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.18"
$ws.Range("E27").Value = "  +9.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.35"
$ws.Range("E28").Value = "  +10.85%  "
$ws.Range("D29").Value = "2.822.34"
$ws.Range("E29").Value = "  +10.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "0.0₃0949"
$ws.Range("E31").Value = "  +11.95%  "
$ws.Range("E32").Value = "  +15.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "518.57"
$ws.Range("E33").Value = "  +14.24%  "
$ws.Range("E34").Value = "  +5.00%  "
$ws.Range("E35").Value = "  +8.92%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.42"
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("E38").Value = "  +7.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.18"
$ws.Range("E39").Value = "  +5.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.37"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.69"
$ws.Range("E44").Value = "  +11.62%  "
$ws.Range("E45").Value = "  +8.80%  "
$ws.Range("E46").Value = "  +10.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.36"
$ws.Range("E47").Value = "  +12.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.13"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0847"
$ws.Range("E49").Value = "  +16.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.65"
$ws.Range("E50").Value = "  +8.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.536"
$ws.Range("E51").Value = "  +9.88%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.00"
$ws.Range("E42").Value = "  +13.91%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "169.98"
$ws.Range("E43").Value = "  +26.99%  "
